$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start clean: remove all existing content/formatting from the old 3x3 table
$ws.Cells.Clear()

# ---- Header row (row 1) ----
# A1 stays blank but carries the header fill/font style (no border)
$headers = @("Cash ratio","Gross profit margin","Return on equity","Operating cash flow","Days of payables outstanding","Ebitgrowth","Interest expense","Current debt","Cost of revenue")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 2).Value = $headers[$i]
}

# ---- Row labels (column A) ----
$tickers = @("TSLA","GM","F","NIO")
for ($i = 0; $i -lt $tickers.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $tickers[$i]
}

# ---- Data values ----
# Cash ratio (B), Gross profit margin (C), Return on equity (D), Operating cash flow (E),
# Days of payables outstanding (F), Ebitgrowth (G), Interest expense (H), Current debt (I), Cost of revenue (J)
$data = @(
    @(0.6085214721629413, 0.25598438535759,    0.2814736936292054,   13956, 91.86878186407959,  0.7622227014997984,  143,     1502,     73825),
    @(0.2100731576234192, 0.1338628895907104,  0.1465364644795846,   22897, 73.90124784536735, -0.5562100638977636,  957,     38778,    150139),
    @(0.2594718477071418, 0.1496928323326395, -0.04976643078488507,  12825, 69.53894060135271, -0.3540207679304516,  9505,    50164,    144310),
    @(0.4337335684024101, 0.1044072101070701, -0.6048686189323729,   $null, 208.6512383532004, -0.9370994324123324,  342.268, 5277.126, 45907.719)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $r + 2
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $col = $c + 2
        $val = $rowVals[$c]
        if ($null -ne $val) {
            $ws.Cells.Item($row, $col).Value = $val
        }
    }
}

# ---- Styling ----

# Header style: bold white font on navy fill
$headerRange = $ws.Range("A1:J1")
$headerRange.Font.Bold = $true
$headerRange.Font.Color = 16777215
$headerRange.Interior.Color = 8388608

# Headers B1:J1 additionally get a thin border and centered alignment
$headerDataRange = $ws.Range("B1:J1")
$headerDataRange.Borders.LineStyle = 1
$headerDataRange.HorizontalAlignment = -4108
$headerDataRange.VerticalAlignment = -4160

# Ticker column (A2:A5): bold font, bordered, centered, percent number format (matches xf 3)
$tickerRange = $ws.Range("A2:A5")
$tickerRange.Font.Bold = $true
$tickerRange.Borders.LineStyle = 1
$tickerRange.HorizontalAlignment = -4108
$tickerRange.VerticalAlignment = -4160
$tickerRange.NumberFormat = "0.00%"

# Percent-formatted data columns: Cash ratio (B), Gross profit margin (C), Days of payables outstanding (F)
$pctRange = $ws.Range("B2:C5")
$pctRange.NumberFormat = "0.00%"
$pctRange2 = $ws.Range("F2:F5")
$pctRange2.NumberFormat = "0.00%"

# NIO's "Operating cash flow" (E5) was missing from the source data (NaN) -- it
# still needs to show up as a present-but-empty cell, matching the known % bug
# that still needs fixing (see commit message).
$ws.Range("E5").Font.Bold = $false

# ---- Column widths (A:I) ----
for ($c = 1; $c -le 9; $c++) {
    $ws.Columns.Item($c).ColumnWidth = 12.17
}
